$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1809.125
$ws.Range("I40").Value = 1735.3529
$ws.Range("J40").Value = 1988.2858
$ws.Range("K40").Value = 1735.3529
$ws.Range("L40").Value = 1988.2858
$ws.Range("M40").Value = -1560.3529
$ws.Range("N40").Value = -2338.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28517.238
$ws.Range("I32").Value = 4894.5747
$ws.Range("K32").Value = 4894.5747
$ws.Range("M32").Value = -4607.5747
$ws.Range("H61").Value = 1844.6938
$ws.Range("I61").Value = 1213.5
$ws.Range("J61").Value = 2450.64
$ws.Range("K61").Value = 1213.5
$ws.Range("L61").Value = 2450.64
$ws.Range("M61").Value = -1001.5
$ws.Range("N61").Value = -2874.64
$ws.Range("H74").Value = 1629.7931
$ws.Range("I74").Value = 1057.7142
$ws.Range("J74").Value = 1811.8182
$ws.Range("K74").Value = 1057.7142
$ws.Range("L74").Value = 1811.8182
$ws.Range("M74").Value = -183.7141999999999
$ws.Range("N74").Value = -3559.8182
$ws.Range("H77").Value = 1629.7931
$ws.Range("I77").Value = 1057.7142
$ws.Range("J77").Value = 1811.8182
$ws.Range("K77").Value = 5288.571
$ws.Range("L77").Value = 9059.091
$ws.Range("M77").Value = -920.5709999999999
$ws.Range("N77").Value = -17795.091
$ws.Range("H110").Value = 23862968
$ws.Range("I110").Value = 29477292
$ws.Range("J110").Value = 2088.5
$ws.Range("K110").Value = 29477292
$ws.Range("L110").Value = 2088.5
$ws.Range("M110").Value = -29475247
$ws.Range("N110").Value = -6178.5
$ws.Range("H132").Value = 3319.647
$ws.Range("I132").Value = 3387.8928
$ws.Range("J132").Value = 3001.1667
$ws.Range("K132").Value = 10163.6784
$ws.Range("L132").Value = 9003.500100000001
$ws.Range("M132").Value = -7633.678400000001
$ws.Range("N132").Value = -14063.5001
$ws.Range("H136").Value = 1844.6938
$ws.Range("I136").Value = 1213.5
$ws.Range("J136").Value = 2450.64
$ws.Range("K136").Value = 3640.5
$ws.Range("L136").Value = 7351.92
$ws.Range("M136").Value = -1090.5
$ws.Range("N136").Value = -12451.92
$ws.Range("H139").Value = 62281.25
$ws.Range("J139").Value = 62281.25
$ws.Range("L139").Value = 62281.25
$ws.Range("N139").Value = -72561.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 252217.62
$ws.Range("I105").Value = 201746
$ws.Range("J105").Value = 336337
$ws.Range("K105").Value = 201746
$ws.Range("L105").Value = 336337
$ws.Range("M105").Value = -199999
$ws.Range("N105").Value = -339831
$ws.Range("H134").Value = 1890.4
$ws.Range("I134").Value = 1984.3334
$ws.Range("J134").Value = 1514.6666
$ws.Range("K134").Value = 5953.0002
$ws.Range("L134").Value = 4543.9998
$ws.Range("M134").Value = -3418.0002
$ws.Range("N134").Value = -9613.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15522.155
$ws.Range("I31").Value = 25323.488
$ws.Range("J31").Value = 2127
$ws.Range("K31").Value = 25323.488
$ws.Range("L31").Value = 2127
$ws.Range("M31").Value = -25028.488
$ws.Range("N31").Value = -2717
$ws.Range("H34").Value = 15522.155
$ws.Range("I34").Value = 25323.488
$ws.Range("J34").Value = 2127
$ws.Range("K34").Value = 25323.488
$ws.Range("L34").Value = 2127
$ws.Range("M34").Value = -25121.488
$ws.Range("N34").Value = -2531
$ws.Range("H86").Value = 2256.6
$ws.Range("I86").Value = 1673.6875
$ws.Range("K86").Value = 1673.6875
$ws.Range("M86").Value = -550.6875
$ws.Range("H89").Value = 2256.6
$ws.Range("I89").Value = 1673.6875
$ws.Range("K89").Value = 8368.4375
$ws.Range("M89").Value = -2752.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 17819.098
$ws.Range("J68").Value = 22244.775
$ws.Range("L68").Value = 66734.32500000001
$ws.Range("N68").Value = -68356.32500000001
$ws.Range("H71").Value = 17819.098
$ws.Range("J71").Value = 22244.775
$ws.Range("L71").Value = 200202.975
$ws.Range("N71").Value = -208314.975
$ws.Range("H86").Value = 37501130
$ws.Range("I86").Value = 1125
$ws.Range("J86").Value = 56251130
$ws.Range("K86").Value = 3375
$ws.Range("L86").Value = 168753390
$ws.Range("M86").Value = -2189
$ws.Range("N86").Value = -168755762
$ws.Range("H89").Value = 37501130
$ws.Range("I89").Value = 1125
$ws.Range("J89").Value = 56251130
$ws.Range("K89").Value = 10125
$ws.Range("L89").Value = 506260170
$ws.Range("M89").Value = -4197
$ws.Range("N89").Value = -506272026
$ws.Range("H131").Value = 2264.9
$ws.Range("J131").Value = 2448.2222
$ws.Range("L131").Value = 7344.6666
$ws.Range("N131").Value = -17424.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 36371.25
$ws.Range("J64").Value = 36371.25
$ws.Range("L64").Value = 36371.25
$ws.Range("N64").Value = -36867.25
$ws.Range("H67").Value = 36371.25
$ws.Range("J67").Value = 36371.25
$ws.Range("L67").Value = 36371.25
$ws.Range("N67").Value = -38087.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2012.3572
$ws.Range("I61").Value = 1743.2858
$ws.Range("J61").Value = 2281.4285
$ws.Range("K61").Value = 1743.2858
$ws.Range("L61").Value = 2281.4285
$ws.Range("M61").Value = -1541.2858
$ws.Range("N61").Value = -2685.4285
$ws.Range("H113").Value = 2012.3572
$ws.Range("I113").Value = 1743.2858
$ws.Range("J113").Value = 2281.4285
$ws.Range("K113").Value = 1743.2858
$ws.Range("L113").Value = 2281.4285
$ws.Range("M113").Value = 426.7141999999999
$ws.Range("N113").Value = -6621.4285
$ws.Range("H122").Value = 3501
$ws.Range("I122").Value = 3501
$ws.Range("K122").Value = 10503
$ws.Range("M122").Value = -8053
$ws.Range("H132").Value = 7265.684
$ws.Range("I132").Value = 8420.916999999999
$ws.Range("J132").Value = 5285.2856
$ws.Range("K132").Value = 25262.751
$ws.Range("L132").Value = 15855.8568
$ws.Range("M132").Value = -22732.751
$ws.Range("N132").Value = -20915.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 33320.65
$ws.Range("I107").Value = 10494.3
$ws.Range("K107").Value = 31482.9
$ws.Range("M107").Value = -29562.9
$ws.Range("H132").Value = 2985.9062
$ws.Range("I132").Value = 3879.8235
$ws.Range("J132").Value = 1972.8
$ws.Range("K132").Value = 11639.4705
$ws.Range("L132").Value = 5918.4
$ws.Range("M132").Value = -10978.4
$ws.Range("N132").Value = -10978.4
